# "Gestione: aggiornamento file ore"
# Fix the typo "GEstione" -> "Gestione" (this merges back into the existing
# shared string already used elsewhere) and append the newly logged hours
# entries to the "Prot. 2.0" sheet.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Prot. 2.0")

# Fix the capitalization typo so it matches the existing "Gestione" entries
$ws2.Range("B5").Value = "Gestione"

# Bring down the formatting (date / text / hours number formats, borders...)
# of the last filled entry so the new rows look consistent with the rest.
$ws2.Range("A5:D6").Copy()
$ws2.Range("A7").PasteSpecial(-4104)
$excel.CutCopyMode = 0

# New log entries
$ws2.Range("A7").Value = 45671
$ws2.Range("B7").Value = "Gestione"
$ws2.Range("C7").Value = "Riflessione su pillars, obiettivi, scopo delle varie arre di gioco"
$ws2.Range("D7").Value = 0.083333333333333329

$ws2.Range("C8").Value = "Riordino Miro e riflessione su palette, contenuti, etc."
$ws2.Range("D8").Value = 0.052083333333333336

$ws2.Activate()
$ws2.Range("D13").Select()
